$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - new data row (Calvin / NoKlein)
$ws.Range("A7").Formula = "=ROW(A6)"
$ws.Range("B7").Value = "Calvin"
$ws.Range("C7").Value = "NoKlein"

# Row 8 - new data row (Anoder / Naame)
$ws.Range("A8").Formula = "=ROW(A7)"
$ws.Range("B8").Value = "Anoder"
$ws.Range("C8").Value = "Naame"

# Row 6's formula becomes its own (non-shared) formula
$ws.Range("A6").Formula = "=ROW(A5)"

# Move the selection to match the author's final cursor position
$ws.Range("J17").Select() | Out-Null
